# Added mapChgWeight into weightChgSymbol.ConfigVal
# The underlying data rows in the worksheet got re-ordered; apply the new
# per-row values (columns A:F) for every row whose contents changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @{
    2  = @(601, 9, 60, 67, 60, 42)
    3  = @(1202, 2, 10, 10, 10, 10)
    4  = @(301, 6, 45, 30, 60, 45)
    5  = @(902, 1, 0, 0, 0, 0)
    6  = @(501, 9, 52, 30, 75, 45)
    7  = @(1203, 3, 15, 15, 15, 15)
    8  = @(101, 9, 30, 15, 60, 15)
    9  = @(701, 3, 90, 45, 97, 15)
    10 = @(201, 9, 30, 15, 45, 30)
    13 = @(401, 9, 48, 67, 75, 45)
    14 = @(1001, 18, 30, 75, 60, 72)
    16 = @(802, 0, 4, 5, 4, 0)
    18 = @(1, 0, 2, 2, 2, 2)
    19 = @(502, 0, 4, 0, 0, 0)
    20 = @(1101, 0, 15, 30, 30, 0)
    21 = @(3, 0, 3, 3, 3, 3)
    22 = @(602, 0, 0, 4, 0, 9)
    23 = @(402, 0, 0, 4, 0, 0)
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    for ($col = 1; $col -le $vals.Length; $col++) {
        $ws.Cells.Item($r, $col).Value = $vals[$col - 1]
    }
}
